$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 434.55554
$ws.Range("I5").Value = 433.64706
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 433.64706
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = -318.64706
$ws.Range("N5").Value = -680
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H137").Value = 1821026.4
$ws.Range("I137").Value = 1690.5555
$ws.Range("K137").Value = 5071.666499999999
$ws.Range("M137").Value = -2521.666499999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5973.1816
$ws.Range("I32").Value = 3519.2856
$ws.Range("J32").Value = 12516.904
$ws.Range("K32").Value = 3519.2856
$ws.Range("L32").Value = 12516.904
$ws.Range("M32").Value = -3232.2856
$ws.Range("N32").Value = -13090.904
$ws.Range("H61").Value = 15827.739
$ws.Range("I61").Value = 27002.285
$ws.Range("K61").Value = 27002.285
$ws.Range("M61").Value = -26790.285
$ws.Range("H63").Value = 6507.294
$ws.Range("I63").Value = 3266
$ws.Range("J63").Value = 7201.857
$ws.Range("K63").Value = 3266
$ws.Range("L63").Value = 7201.857
$ws.Range("M63").Value = -2580
$ws.Range("N63").Value = -8573.857
$ws.Range("H66").Value = 6507.294
$ws.Range("I66").Value = 3266
$ws.Range("J66").Value = 7201.857
$ws.Range("K66").Value = 16330
$ws.Range("L66").Value = 36009.285
$ws.Range("M66").Value = -12898
$ws.Range("N66").Value = -42873.285
$ws.Range("H74").Value = 2671.9822
$ws.Range("I74").Value = 2309.925
$ws.Range("J74").Value = 3577.125
$ws.Range("K74").Value = 2309.925
$ws.Range("L74").Value = 3577.125
$ws.Range("M74").Value = -1435.925
$ws.Range("N74").Value = -5325.125
$ws.Range("H77").Value = 2671.9822
$ws.Range("I77").Value = 2309.925
$ws.Range("J77").Value = 3577.125
$ws.Range("K77").Value = 11549.625
$ws.Range("L77").Value = 17885.625
$ws.Range("M77").Value = -7181.625
$ws.Range("N77").Value = -26621.625
$ws.Range("H132").Value = 4894.778
$ws.Range("I132").Value = 3800.5
$ws.Range("J132").Value = 8724.75
$ws.Range("K132").Value = 11401.5
$ws.Range("L132").Value = 26174.25
$ws.Range("M132").Value = -8871.5
$ws.Range("N132").Value = -31234.25
$ws.Range("H136").Value = 15827.739
$ws.Range("I136").Value = 27002.285
$ws.Range("K136").Value = 81006.855
$ws.Range("M136").Value = -78456.855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 95135.75
$ws.Range("J62").Value = 95135.75
$ws.Range("L62").Value = 95135.75
$ws.Range("N62").Value = -96507.75
$ws.Range("H65").Value = 95135.75
$ws.Range("J65").Value = 95135.75
$ws.Range("L65").Value = 285407.25
$ws.Range("N65").Value = -292271.25
$ws.Range("H86").Value = 1275.35
$ws.Range("J86").Value = 1072.875
$ws.Range("L86").Value = 1072.875
$ws.Range("N86").Value = -3318.875
$ws.Range("H89").Value = 1275.35
$ws.Range("J89").Value = 1072.875
$ws.Range("L89").Value = 5364.375
$ws.Range("N89").Value = -16596.375
$ws.Range("H107").Value = 1440.8214
$ws.Range("I107").Value = 1509.2
$ws.Range("K107").Value = 1509.2
$ws.Range("M107").Value = 410.8
$ws.Range("H124").Value = 43985
$ws.Range("J124").Value = 43985
$ws.Range("L124").Value = 43985
$ws.Range("N124").Value = -53805

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 25878.5
$ws.Range("I35").Value = 23000
$ws.Range("J35").Value = 28757
$ws.Range("K35").Value = 23000
$ws.Range("L35").Value = 28757
$ws.Range("M35").Value = -22706
$ws.Range("N35").Value = -29345
$ws.Range("H99").Value = 5953.25
$ws.Range("I99").Value = 5511.4287
$ws.Range("J99").Value = 6571.8
$ws.Range("K99").Value = 5511.4287
$ws.Range("L99").Value = 6571.8
$ws.Range("M99").Value = -4013.4287
$ws.Range("N99").Value = -9567.799999999999
$ws.Range("H126").Value = 5953.25
$ws.Range("I126").Value = 5511.4287
$ws.Range("J126").Value = 6571.8
$ws.Range("K126").Value = 16534.2861
$ws.Range("L126").Value = 19715.4
$ws.Range("M126").Value = -14064.2861
$ws.Range("N126").Value = -24655.4
$ws.Range("H134").Value = 2238.0588
$ws.Range("I134").Value = 2331.9285
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 6995.7855
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -4460.7855
$ws.Range("N134").Value = -10470

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3970.5715
$ws.Range("I70").Value = 1834.4445
$ws.Range("K70").Value = 5503.333500000001
$ws.Range("M70").Value = -5188.333500000001
$ws.Range("H73").Value = 3970.5715
$ws.Range("I73").Value = 1834.4445
$ws.Range("K73").Value = 5503.333500000001
$ws.Range("M73").Value = -4411.333500000001
$ws.Range("H75").Value = 4239.0713
$ws.Range("I75").Value = 1311.6666
$ws.Range("J75").Value = 6434.625
$ws.Range("K75").Value = 3934.9998
$ws.Range("L75").Value = 19303.875
$ws.Range("M75").Value = -2936.9998
$ws.Range("N75").Value = -21299.875
$ws.Range("H78").Value = 4239.0713
$ws.Range("I78").Value = 1311.6666
$ws.Range("J78").Value = 6434.625
$ws.Range("K78").Value = 11804.9994
$ws.Range("L78").Value = 57911.625
$ws.Range("M78").Value = -6812.999400000001
$ws.Range("N78").Value = -67895.625
$ws.Range("H87").Value = 19001.1
$ws.Range("I87").Value = 14991.167
$ws.Range("K87").Value = 44973.501
$ws.Range("M87").Value = -43725.501
$ws.Range("H90").Value = 19001.1
$ws.Range("I90").Value = 14991.167
$ws.Range("K90").Value = 134920.503
$ws.Range("M90").Value = -128680.503
$ws.Range("H103").Value = 216.42857
$ws.Range("J103").Value = 64
$ws.Range("L103").Value = 192
$ws.Range("N103").Value = -1950
$ws.Range("H121").Value = 9524574
$ws.Range("I121").Value = 20000364
$ws.Range("J121").Value = 1128.4546
$ws.Range("K121").Value = 60001092
$ws.Range("L121").Value = 3385.3638
$ws.Range("M121").Value = -59999782
$ws.Range("N121").Value = -6005.3638
$ws.Range("H132").Value = 2824.3157
$ws.Range("I132").Value = 2154.5715
$ws.Range("K132").Value = 19391.1435
$ws.Range("M132").Value = -16861.1435

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 79531.42999999999
$ws.Range("J95").Value = 79531.42999999999
$ws.Range("L95").Value = 79531.42999999999
$ws.Range("N95").Value = -85023.42999999999
$ws.Range("H132").Value = 672786.5600000001
$ws.Range("I132").Value = 3013173
$ws.Range("J132").Value = 4104.7144
$ws.Range("K132").Value = 9039519
$ws.Range("L132").Value = 12314.1432
$ws.Range("M132").Value = -9036989
$ws.Range("N132").Value = -17374.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3341.3794
$ws.Range("J46").Value = 3403.7036
$ws.Range("L46").Value = 3403.7036
$ws.Range("N46").Value = -3779.7036
$ws.Range("H136").Value = 6568.9443
$ws.Range("I136").Value = 6279.5835
$ws.Range("J136").Value = 7147.6665
$ws.Range("K136").Value = 18838.7505
$ws.Range("L136").Value = 21442.9995
$ws.Range("M136").Value = -16288.7505
$ws.Range("N136").Value = -26542.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 71607
$ws.Range("J46").Value = 71607
$ws.Range("L46").Value = 71607
$ws.Range("N46").Value = -72069
$ws.Range("H55").Value = 10499
$ws.Range("I55").Value = 7998
$ws.Range("K55").Value = 7998
$ws.Range("M55").Value = -7721
$ws.Range("H81").Value = 645.1111
$ws.Range("I81").Value = 566.5714
$ws.Range("J81").Value = 920
$ws.Range("K81").Value = 1133.1428
$ws.Range("L81").Value = 1840
$ws.Range("M81").Value = -72.14280000000008
$ws.Range("N81").Value = -3962
$ws.Range("H84").Value = 645.1111
$ws.Range("I84").Value = 566.5714
$ws.Range("J84").Value = 920
$ws.Range("K84").Value = 5665.714
$ws.Range("L84").Value = 9200
$ws.Range("M84").Value = -361.7139999999999
$ws.Range("N84").Value = -19808
$ws.Range("H134").Value = 71607
$ws.Range("J134").Value = 71607
$ws.Range("L134").Value = 214821
$ws.Range("N134").Value = -219891
$ws.Range("H136").Value = 9774495
$ws.Range("J136").Value = 3828
$ws.Range("L136").Value = 11484
$ws.Range("N136").Value = -16584
